# feat: introduce state transition dates
# Replace the ingest_session_issue_id UUID values in column A (rows 2-49)
# of the "diagnostics" worksheet with new UUID values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diagnostics")

$newIds = @(
    "925f93cb-4100-4fa4-931f-def238dacd7d",
    "6f3cf552-3f63-446d-be3a-04553cc2865e",
    "3f52318a-fd2c-45c3-b3d3-bbcb6a363b87",
    "625c85d6-6eee-4711-a935-b1902bbb42b6",
    "55d6ce71-541d-4b50-b95b-82657a450e67",
    "f1195c33-36cd-4bed-97ae-0f202c817eb1",
    "35f78a04-30f6-4846-98d7-738f539c5a45",
    "98808950-8749-4246-b2ed-81547dbb0546",
    "7425f039-ddf7-421b-be64-551f6fe373f8",
    "893ef79d-048e-48d3-97ec-7369791b7c94",
    "5875a951-5b4f-4d9e-b063-22f5668dadc8",
    "92d7eb46-7dac-424a-9523-17d8c3dc14f9",
    "b1ec777e-169b-4451-8ee2-db15d49bb32a",
    "56ff0086-bd1b-463c-88f0-a0e5186d8035",
    "ce5d302d-74ee-4798-ba28-4a4c7d2d30f8",
    "d7cc8d6b-ab6a-490b-acb2-9c395128ec10",
    "259f86ec-768b-4496-bfee-0e6f912bef6d",
    "12adb95b-e6eb-452d-8c52-2137bf952663",
    "8a7b63e6-4ccb-4532-8e93-2cd968f2060d",
    "02983050-cc92-49cf-98d0-97fb5b933a47",
    "b2d996fa-f538-416e-9425-cc5d811db416",
    "54776c82-8bd0-46f6-b5df-a6a3fae3e156",
    "98708eca-7f94-4718-af91-6a08eb96e5f3",
    "571a8385-db8f-404c-981b-1829c4254a2e",
    "0e65a070-6f89-4b6d-a477-70a0169832ef",
    "51595ce7-95d9-40cb-99aa-4bf606b84974",
    "fe282f1b-a01d-42de-9bc2-523a36d54d78",
    "5b1fc05f-35b4-4409-bc1f-87a8f583d26e",
    "a1ed2d29-ccba-4874-ba3c-bdcc737cfc4d",
    "626781d0-1f85-41a0-84f1-8b27b014199b",
    "8dc4a1fc-bf16-414b-b771-989fb0a1eb79",
    "8008a107-bc83-4156-8647-230c07a66f8f",
    "c0a3b1c9-37d0-4623-87be-17fa0be638a9",
    "0114fec6-986c-4b95-8b75-17e3abf4e273",
    "d3292789-430a-43dc-8836-155a4d52c4c6",
    "12f52213-e5e4-42e3-8033-3b2aecaac1ae",
    "8cd657c5-7d3e-474e-855e-04075ffa5196",
    "fc915eb0-6cdd-45c1-92e7-982509ab6225",
    "4ccc6510-50a1-4605-b7e2-ee04f599f8e7",
    "43848a3e-6c56-4cda-a4f1-769573162373",
    "8a465c88-30d3-4d4c-9988-8ab80c164701",
    "38c871bc-f811-4f6b-a2c3-1cb59602ceb9",
    "1e6844d2-b4a2-447e-8220-0ed7af21833d",
    "22098fea-77de-4093-a4c7-ddc00a27ca2f",
    "f89aaa40-3e1d-4989-858f-026dacb64532",
    "353c9d9a-7acd-4e91-acca-a02081b47437",
    "13d8c905-b8bd-4e9a-909b-ba8eb64e7c4f",
    "64d47455-29fd-4d20-acc3-37ce0e208b1f"
)

$startRow = 2
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newIds[$i]
}
